$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.789999999999999
$ws.Range("D6").Value = -8.169
$ws.Range("D7").Value = -7.842000000000001
$ws.Range("D8").Value = -8.323
$ws.Range("D16").Value = -8.529
$ws.Range("D20").Value = -7.797
$ws.Range("D21").Value = -8.076000000000001
$ws.Range("D28").Value = -7.945
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.243
$ws.Range("D32").Value = -7.665000000000001
$ws.Range("D40").Value = -8.119000000000002
$ws.Range("D46").Value = -8.202000000000002
$ws.Range("D51").Value = -8.457000000000001
$ws.Range("D52").Value = -8.115
$ws.Range("D57").Value = -8.06
$ws.Range("D59").Value = -8.123999999999999
$ws.Range("D62").Value = -7.946
$ws.Range("D66").Value = -7.061000000000002
$ws.Range("D73").Value = -7.902999999999999
$ws.Range("D74").Value = -7.911
$ws.Range("D77").Value = -7.840000000000001
$ws.Range("D92").Value = -6.536999999999999
$ws.Range("D100").Value = -8.295
